$wb = $excel.ActiveWorkbook

# Sheet "Bank and Cash Accounts": update Txn Date Min (column F) values
$wsBank = $wb.Worksheets.Item("Bank and Cash Accounts")
$wsBank.Range("F2").Value = 43104
$wsBank.Range("F3").Value = 43101

# Sheet "Credit Cards": update Txn Date Min (column F) values
$wsCredit = $wb.Worksheets.Item("Credit Cards")
$wsCredit.Range("F2").Value = 43101
$wsCredit.Range("F3").Value = 43102
